# Add a new weekly price record for "Acelga" (Feria Lagunitas de Puerto Montt)
# dated 2023-10-05. The new record is inserted as row 286, pushing every
# existing record from row 286 onward down by one row (old row 286 becomes
# 287, ..., old row 354 becomes 355).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 286, shifting rows 286:354 down to 287:355.
$ws.Rows.Item(286).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(286, 1).Value  = 4
$ws.Cells.Item(286, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(286, 3).Value  = "Los Lagos"
$ws.Cells.Item(286, 4).Value  = "2023-10-05"
$ws.Cells.Item(286, 5).Value  = 10
$ws.Cells.Item(286, 6).Value  = 100112009
$ws.Cells.Item(286, 7).Value  = "Acelga"
$ws.Cells.Item(286, 8).Value  = "Sin especificar"
$ws.Cells.Item(286, 9).Value  = "Primera"
$ws.Cells.Item(286, 10).Value = 75
$ws.Cells.Item(286, 11).Value = 10000
$ws.Cells.Item(286, 12).Value = 10000
$ws.Cells.Item(286, 13).Value = 10000
$ws.Cells.Item(286, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(286, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(286, 16).Value = 833
$ws.Cells.Item(286, 17).Value = 12
$ws.Cells.Item(286, 18).Value = "Hortaliza"
